$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 11198
$ws.Range("I11").Value = 11198
$ws.Range("K11").Value = 11198
$ws.Range("M11").Value = -11058
$ws.Range("H17").Value = 5642613
$ws.Range("J17").Value = 5642613
$ws.Range("L17").Value = 16927839
$ws.Range("N17").Value = -16928175
$ws.Range("H33").Value = 13901514
$ws.Range("I33").Value = 16981
$ws.Range("J33").Value = 41670580
$ws.Range("K33").Value = 16981
$ws.Range("L33").Value = 41670580
$ws.Range("M33").Value = -16752
$ws.Range("N33").Value = -41671038
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2864
$ws.Range("H61").Value = 633
$ws.Range("I61").Value = 499.5
$ws.Range("K61").Value = 1498.5
$ws.Range("M61").Value = -1326.5
$ws.Range("H64").Value = 6301.5264
$ws.Range("I64").Value = 7061.125
$ws.Range("J64").Value = 5749.091
$ws.Range("K64").Value = 7061.125
$ws.Range("L64").Value = 5749.091
$ws.Range("M64").Value = -6813.125
$ws.Range("N64").Value = -6245.091
$ws.Range("H67").Value = 6301.5264
$ws.Range("I67").Value = 7061.125
$ws.Range("J67").Value = 5749.091
$ws.Range("K67").Value = 7061.125
$ws.Range("L67").Value = 5749.091
$ws.Range("M67").Value = -6203.125
$ws.Range("N67").Value = -7465.091
$ws.Range("H70").Value = 1891.6666
$ws.Range("I70").Value = 1337.5
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 4012.5
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -3742.5
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 1891.6666
$ws.Range("I73").Value = 1337.5
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 4012.5
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3076.5
$ws.Range("N73").Value = -10872
$ws.Range("H138").Value = 6122.9766
$ws.Range("I138").Value = 17487
$ws.Range("J138").Value = 3114.853
$ws.Range("K138").Value = 52461
$ws.Range("L138").Value = 9344.559000000001
$ws.Range("M138").Value = -47321
$ws.Range("N138").Value = -19624.559

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 810
$ws.Range("I2").Value = 777.1429000000001
$ws.Range("K2").Value = 777.1429000000001
$ws.Range("M2").Value = -664.1429000000001
$ws.Range("H61").Value = 1843630.2
$ws.Range("I61").Value = 5569.032
$ws.Range("K61").Value = 5569.032
$ws.Range("M61").Value = -5357.032
$ws.Range("H63").Value = 7470.9697
$ws.Range("I63").Value = 3317.6296
$ws.Range("K63").Value = 3317.6296
$ws.Range("M63").Value = -2631.6296
$ws.Range("H66").Value = 7470.9697
$ws.Range("I66").Value = 3317.6296
$ws.Range("K66").Value = 16588.148
$ws.Range("M66").Value = -13156.148
$ws.Range("H74").Value = 477414.84
$ws.Range("I74").Value = 3397.44
$ws.Range("K74").Value = 3397.44
$ws.Range("M74").Value = -2523.44
$ws.Range("H77").Value = 477414.84
$ws.Range("I77").Value = 3397.44
$ws.Range("K77").Value = 16987.2
$ws.Range("M77").Value = -12619.2
$ws.Range("H109").Value = 70149.7
$ws.Range("J109").Value = 70149.7
$ws.Range("L109").Value = 70149.7
$ws.Range("N109").Value = -72923.7
$ws.Range("H116").Value = 810
$ws.Range("I116").Value = 777.1429000000001
$ws.Range("K116").Value = 777.1429000000001
$ws.Range("M116").Value = 1516.8571
$ws.Range("H131").Value = 86333.336
$ws.Range("J131").Value = 86333.336
$ws.Range("L131").Value = 86333.336
$ws.Range("N131").Value = -96413.336
$ws.Range("H136").Value = 1843630.2
$ws.Range("I136").Value = 5569.032
$ws.Range("K136").Value = 16707.096
$ws.Range("M136").Value = -14157.096

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 810
$ws.Range("I3").Value = 777.1429000000001
$ws.Range("K3").Value = 777.1429000000001
$ws.Range("M3").Value = -663.1429000000001
$ws.Range("H105").Value = 17640.111
$ws.Range("I105").Value = 29943.25
$ws.Range("K105").Value = 29943.25
$ws.Range("M105").Value = -28196.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2832.606
$ws.Range("I31").Value = 1708.0588
$ws.Range("J31").Value = 3222.7551
$ws.Range("K31").Value = 1708.0588
$ws.Range("L31").Value = 3222.7551
$ws.Range("M31").Value = -1413.0588
$ws.Range("N31").Value = -3812.7551
$ws.Range("H34").Value = 2832.606
$ws.Range("I34").Value = 1708.0588
$ws.Range("J34").Value = 3222.7551
$ws.Range("K34").Value = 1708.0588
$ws.Range("L34").Value = 3222.7551
$ws.Range("M34").Value = -1506.0588
$ws.Range("N34").Value = -3626.7551
$ws.Range("H94").Value = 1277.421
$ws.Range("I94").Value = 997.5
$ws.Range("K94").Value = 997.5
$ws.Range("M94").Value = -546.5
$ws.Range("H132").Value = 15155655
$ws.Range("I132").Value = 4108.7334
$ws.Range("K132").Value = 12326.2002
$ws.Range("M132").Value = -9796.200199999999
$ws.Range("H141").Value = 701986.25
$ws.Range("J141").Value = 681318.5
$ws.Range("L141").Value = 681318.5
$ws.Range("N141").Value = -691678.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 52632050
$ws.Range("J107").Value = 76923420
$ws.Range("L107").Value = 230770260
$ws.Range("N107").Value = -230774100
$ws.Range("H140").Value = 10871345
$ws.Range("I140").Value = 11906235
$ws.Range("K140").Value = 35718705
$ws.Range("M140").Value = -35713525

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12276.238
$ws.Range("I70").Value = 11083.941
$ws.Range("J70").Value = 17343.5
$ws.Range("K70").Value = 11083.941
$ws.Range("L70").Value = 17343.5
$ws.Range("M70").Value = -10813.941
$ws.Range("N70").Value = -17883.5
$ws.Range("H73").Value = 12276.238
$ws.Range("I73").Value = 11083.941
$ws.Range("J73").Value = 17343.5
$ws.Range("K73").Value = 11083.941
$ws.Range("L73").Value = 17343.5
$ws.Range("M73").Value = -10147.941
$ws.Range("N73").Value = -19215.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H122").Value = 8579.909
$ws.Range("I122").Value = 5658.2856
$ws.Range("J122").Value = 13692.75
$ws.Range("K122").Value = 16974.8568
$ws.Range("L122").Value = 41078.25
$ws.Range("M122").Value = -14524.8568
$ws.Range("N122").Value = -45978.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4695.579
$ws.Range("J22").Value = 5880.846
$ws.Range("L22").Value = 5880.846
$ws.Range("N22").Value = -6470.846
$ws.Range("H27").Value = 4695.579
$ws.Range("J27").Value = 5880.846
$ws.Range("L27").Value = 5880.846
$ws.Range("N27").Value = -6094.846
$ws.Range("H55").Value = 364.68
$ws.Range("I55").Value = 138.76923
$ws.Range("J55").Value = 609.4167
$ws.Range("K55").Value = 138.76923
$ws.Range("L55").Value = 609.4167
$ws.Range("M55").Value = 34.23077000000001
$ws.Range("N55").Value = -955.4167
$ws.Range("H61").Value = 3144.8
$ws.Range("I61").Value = 3182.3914
$ws.Range("J61").Value = 3021.2856
$ws.Range("K61").Value = 3182.3914
$ws.Range("L61").Value = 3021.2856
$ws.Range("M61").Value = -2980.3914
$ws.Range("N61").Value = -3425.2856
$ws.Range("H105").Value = 22179.334
$ws.Range("J105").Value = 22179.334
$ws.Range("L105").Value = 22179.334
$ws.Range("N105").Value = -29167.334
$ws.Range("H113").Value = 3144.8
$ws.Range("I113").Value = 3182.3914
$ws.Range("J113").Value = 3021.2856
$ws.Range("K113").Value = 3182.3914
$ws.Range("L113").Value = 3021.2856
$ws.Range("M113").Value = -1012.3914
$ws.Range("N113").Value = -7361.2856

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 51512.89
$ws.Range("J46").Value = 51512.89
$ws.Range("L46").Value = 51512.89
$ws.Range("N46").Value = -51974.89
$ws.Range("H107").Value = 76923496
$ws.Range("I107").Value = 499.125
$ws.Range("J107").Value = 200000290
$ws.Range("K107").Value = 1497.375
$ws.Range("L107").Value = 600000870
$ws.Range("M107").Value = 422.625
$ws.Range("N107").Value = -600004710
$ws.Range("H132").Value = 2123.7742
$ws.Range("I132").Value = 2055.4736
$ws.Range("J132").Value = 2231.9167
$ws.Range("K132").Value = 6166.4208
$ws.Range("L132").Value = 6695.750100000001
$ws.Range("M132").Value = -3636.4208
$ws.Range("N132").Value = -11755.7501
$ws.Range("H134").Value = 51512.89
$ws.Range("J134").Value = 51512.89
$ws.Range("L134").Value = 154538.67
$ws.Range("N134").Value = -159608.67
$ws.Range("H136").Value = 2349.111
$ws.Range("I136").Value = 1372.1111
$ws.Range("K136").Value = 4116.3333
$ws.Range("M136").Value = -1566.3333

Write-Host "Done applying market price updates"